# Auto-generated edit script: updates Leve profitability calc columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 1139.6
$ws.Range("J2").Value = 1139.6
$ws.Range("L2").Value = 1139.6
$ws.Range("N2").Value = -1365.6
# row 69
$ws.Range("H69").Value = 9962.714
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9962.714
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 29888.142
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -31636.142
# row 72
$ws.Range("H72").Value = 9962.714
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9962.714
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 89664.42600000001
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -98400.42600000001
# row 76
$ws.Range("H76").Value = 4375.6665
$ws.Range("I76").Value = 4156.5713
$ws.Range("K76").Value = 4156.5713
$ws.Range("M76").Value = -3841.5713
# row 79
$ws.Range("H79").Value = 4375.6665
$ws.Range("I79").Value = 4156.5713
$ws.Range("K79").Value = 4156.5713
$ws.Range("M79").Value = -3064.5713
# row 99
$ws.Range("H99").Value = 1908.4
$ws.Range("I99").Value = 2259.4
$ws.Range("K99").Value = 6778.200000000001
$ws.Range("M99").Value = -5280.200000000001
# row 100
$ws.Range("H100").Value = 3696.68
$ws.Range("I100").Value = 1511.8
$ws.Range("J100").Value = 6974
$ws.Range("K100").Value = 1511.8
$ws.Range("L100").Value = 6974
$ws.Range("M100").Value = -970.8
$ws.Range("N100").Value = -8056
# row 125
$ws.Range("H125").Value = 2896.9
$ws.Range("I125").Value = 2459.3333
$ws.Range("K125").Value = 22133.9997
$ws.Range("M125").Value = -19673.9997
# row 135
$ws.Range("H135").Value = 1333.6666
$ws.Range("I135").Value = 1350.3704
$ws.Range("K135").Value = 12153.3336
$ws.Range("M135").Value = -9618.3336
# row 137
$ws.Range("H137").Value = 3049
$ws.Range("I137").Value = 2046.4
$ws.Range("J137").Value = 3638.7646
$ws.Range("K137").Value = 6139.200000000001
$ws.Range("L137").Value = 10916.2938
$ws.Range("M137").Value = -3589.200000000001
$ws.Range("N137").Value = -16016.2938
# row 138
$ws.Range("H138").Value = 3236.2068
$ws.Range("J138").Value = 3604.3948
$ws.Range("L138").Value = 10813.1844
$ws.Range("N138").Value = -21093.1844

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
# row 32
$ws.Range("H32").Value = 5107609.5
$ws.Range("I32").Value = 5955878
$ws.Range("K32").Value = 5955878
$ws.Range("M32").Value = -5955591
# row 61
$ws.Range("H61").Value = 5835.6055
$ws.Range("I61").Value = 3187.8928
$ws.Range("J61").Value = 13249.2
$ws.Range("K61").Value = 3187.8928
$ws.Range("L61").Value = 13249.2
$ws.Range("M61").Value = -2975.8928
$ws.Range("N61").Value = -13673.2
# row 74
$ws.Range("H74").Value = 2722098.2
$ws.Range("I74").Value = 3791748.5
$ws.Range("J74").Value = 6832.5386
$ws.Range("K74").Value = 3791748.5
$ws.Range("L74").Value = 6832.5386
$ws.Range("M74").Value = -3790874.5
$ws.Range("N74").Value = -8580.5386
# row 77
$ws.Range("H77").Value = 2722098.2
$ws.Range("I77").Value = 3791748.5
$ws.Range("J77").Value = 6832.5386
$ws.Range("K77").Value = 18958742.5
$ws.Range("L77").Value = 34162.693
$ws.Range("M77").Value = -18954374.5
$ws.Range("N77").Value = -42898.693
# row 88
$ws.Range("H88").Value = 2474.4285
$ws.Range("J88").Value = 2832.125
$ws.Range("L88").Value = 2832.125
$ws.Range("N88").Value = -3644.125
# row 91
$ws.Range("H91").Value = 2474.4285
$ws.Range("J91").Value = 2832.125
$ws.Range("L91").Value = 2832.125
$ws.Range("N91").Value = -5640.125
# row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
# row 122
$ws.Range("H122").Value = 1837.6471
$ws.Range("I122").Value = 1826.6562
$ws.Range("J122").Value = 2013.5
$ws.Range("K122").Value = 5479.9686
$ws.Range("L122").Value = 6040.5
$ws.Range("M122").Value = -3029.9686
$ws.Range("N122").Value = -10940.5
# row 132
$ws.Range("H132").Value = 867246.2
$ws.Range("I132").Value = 1296185.1
$ws.Range("J132").Value = 9368.333000000001
$ws.Range("K132").Value = 3888555.3
$ws.Range("L132").Value = 28104.999
$ws.Range("M132").Value = -3886025.3
$ws.Range("N132").Value = -33164.999
# row 136
$ws.Range("H136").Value = 5835.6055
$ws.Range("I136").Value = 3187.8928
$ws.Range("J136").Value = 13249.2
$ws.Range("K136").Value = 9563.678400000001
$ws.Range("L136").Value = 39747.60000000001
$ws.Range("M136").Value = -7013.678400000001
$ws.Range("N136").Value = -44847.60000000001

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row 100
$ws.Range("H100").Value = 20649.2
$ws.Range("J100").Value = 20649.2
$ws.Range("L100").Value = 20649.2
$ws.Range("N100").Value = -22813.2

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 6
$ws.Range("H6").Value = 2667.111
$ws.Range("I6").Value = 800.8
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 800.8
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = -687.8
$ws.Range("N6").Value = -5226
# row 31
$ws.Range("H31").Value = 9021.718000000001
$ws.Range("I31").Value = 1909.8182
$ws.Range("J31").Value = 11815.679
$ws.Range("K31").Value = 1909.8182
$ws.Range("L31").Value = 11815.679
$ws.Range("M31").Value = -1614.8182
$ws.Range("N31").Value = -12405.679
# row 34
$ws.Range("H34").Value = 9021.718000000001
$ws.Range("I34").Value = 1909.8182
$ws.Range("J34").Value = 11815.679
$ws.Range("K34").Value = 1909.8182
$ws.Range("L34").Value = 11815.679
$ws.Range("M34").Value = -1707.8182
$ws.Range("N34").Value = -12219.679
# row 99
$ws.Range("H99").Value = 4690.778
$ws.Range("I99").Value = 3002.2
$ws.Range("J99").Value = 6801.5
$ws.Range("K99").Value = 3002.2
$ws.Range("L99").Value = 6801.5
$ws.Range("M99").Value = -1504.2
$ws.Range("N99").Value = -9797.5
# row 126
$ws.Range("H126").Value = 4690.778
$ws.Range("I126").Value = 3002.2
$ws.Range("J126").Value = 6801.5
$ws.Range("K126").Value = 9006.599999999999
$ws.Range("L126").Value = 20404.5
$ws.Range("M126").Value = -6536.599999999999
$ws.Range("N126").Value = -25344.5
# row 132
$ws.Range("H132").Value = 6590879.5
$ws.Range("J132").Value = 35721850
$ws.Range("L132").Value = 107165550
$ws.Range("N132").Value = -107170610
# row 134
$ws.Range("H134").Value = 13285.609
$ws.Range("I134").Value = 11667.389
$ws.Range("J134").Value = 24936.8
$ws.Range("K134").Value = 35002.167
$ws.Range("L134").Value = 74810.39999999999
$ws.Range("M134").Value = -32467.167
$ws.Range("N134").Value = -79880.39999999999

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 28
$ws.Range("H28").Value = 1646
$ws.Range("I28").Value = 1646
$ws.Range("K28").Value = 4938
$ws.Range("M28").Value = -4706

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 5439.1665
$ws.Range("I132").Value = 3320.5715
$ws.Range("K132").Value = 9961.7145
$ws.Range("M132").Value = -7431.7145

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3255.3044
$ws.Range("I7").Value = 3265.3333
$ws.Range("K7").Value = 3265.3333
$ws.Range("M7").Value = -3153.3333
# row 40
$ws.Range("H40").Value = 4446.6216
$ws.Range("I40").Value = 4701.6
$ws.Range("K40").Value = 4701.6
$ws.Range("M40").Value = -4565.6
# row 82
$ws.Range("H82").Value = 2964.5334
$ws.Range("I82").Value = 2725.3635
$ws.Range("K82").Value = 2725.3635
$ws.Range("M82").Value = -2364.3635
# row 85
$ws.Range("H85").Value = 2964.5334
$ws.Range("I85").Value = 2725.3635
$ws.Range("K85").Value = 2725.3635
$ws.Range("M85").Value = -1477.3635
# row 126
$ws.Range("H126").Value = 3255.3044
$ws.Range("I126").Value = 3265.3333
$ws.Range("K126").Value = 9795.999899999999
$ws.Range("M126").Value = -7325.999899999999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 4157.2856
$ws.Range("I62").Value = 4025.5
$ws.Range("K62").Value = 4025.5
$ws.Range("M62").Value = -3401.5
# row 65
$ws.Range("H65").Value = 4157.2856
$ws.Range("I65").Value = 4025.5
$ws.Range("K65").Value = 20127.5
$ws.Range("M65").Value = -17007.5
